$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (MMST 470): quantity cleared
$ws.Range("C9").Value = $null

# Row 31 (SC Data -10): quantity updated
$ws.Range("C31").Value = 7873

# Row 32 (SC Data -29): quantity added
$ws.Range("C32").Value = 500

# Row 43 (I'top-Up): quantity cleared
$ws.Range("C43").Value = $null
